# Applies updated loading_percent results for the 380 kV case (Case_2_235).
# Only columns B, C, E, F, G, H, K, M change for rows 2-25 (data rows 0-23);
# columns D, I, J, L, N, O remain 0 and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$columns = @("B", "C", "E", "F", "G", "H", "K", "M")

# Each entry: worksheet row number, followed by the new values for
# columns B, C, E, F, G, H, K, M (in that order).
$newValues = @(
    @(2, 9.53738884769195, 4.68945060299725, 12.41774513761585, 16.86991607391233, 39.92623383480575, 16.89167979842016, 8.886016925231822, 13.87640065454884),
    @(3, 9.316102732180932, 4.549524839788759, 12.18783809470999, 15.89584955866808, 39.76252730385587, 16.91321953295695, 8.751113932555404, 13.73215210226653),
    @(4, 9.180255287229997, 4.460299028150407, 12.04856273325771, 15.26997757108491, 39.67275695079184, 16.9293109631106, 8.669595702558155, 13.64663779154895),
    @(5, 9.124992962736405, 4.423133520943233, 11.99236391437331, 15.00819731993403, 39.63889391000607, 16.936587196589, 8.636751575998995, 13.61259620671732),
    @(6, 9.115825233556839, 4.4169144292458, 11.98306817853113, 14.96433081551593, 39.63343561974781, 16.93783876986607, 8.631321805010584, 13.6069933963253),
    @(7, 9.179509487187907, 4.459801023635081, 12.0478024527093, 15.26647399323137, 39.67228923322031, 16.92940618493286, 8.669151176702965, 13.64617538166483),
    @(8, 9.461138226889853, 4.641906912236045, 12.33813057152248, 16.53996406344768, 39.86756882039805, 16.89851081919215, 8.83925708840046, 13.82605303045343),
    @(9, 10.00975744400161, 4.971680498030285, 12.91884490214601, 19.00274580682531, 40.33485063257805, 16.86074163129512, 9.181257142790841, 14.20128490995503),
    @(10, 10.40572735139212, 5.196049733866221, 13.34768249034368, 20.67494806633232, 40.72807968292886, 16.84699959967402, 9.435003735403932, 14.48825503970385),
    @(11, 10.58338024607312, 5.294012469400585, 13.54227917689735, 21.3917225636224, 40.91741149312346, 16.84380697772851, 9.550428330598852, 14.62072465115524),
    @(12, 10.65022503058942, 5.330502701108642, 13.61582158028465, 21.65686569030329, 40.99057022754115, 16.84303893759775, 9.594091476460477, 14.67111938141751),
    @(13, 10.6358489669225, 5.322671068955782, 13.59999064934827, 21.60004134736742, 40.97474979042556, 16.84318472254318, 9.584690563201056, 14.6602564075743),
    @(14, 10.58888859133636, 5.297026775269904, 13.54833293532977, 21.4136618050453, 40.92340122714081, 16.84373494796973, 9.554021689631002, 14.62486627922643),
    @(15, 10.56006613841186, 5.281239527057072, 13.51666968987863, 21.29868154950795, 40.89213802453119, 16.84412942922229, 9.53522887220875, 14.60321758570833),
    @(16, 10.39406118119909, 5.189563646018649, 13.33494869450775, 20.62722412089977, 40.71591335507981, 16.8472698922083, 9.427456483474991, 14.47963269111239),
    @(17, 10.29153517040993, 5.132260349037203, 13.22328857302425, 20.20408069597325, 40.6104544294139, 16.84998071346143, 9.361307744734532, 14.40427811207085),
    @(18, 10.23233589369485, 5.098915552290003, 13.15902353208303, 19.95656407809801, 40.55078343926385, 16.85182775132982, 9.323262766150956, 14.36112011217225),
    @(19, 10.21225499458289, 5.08755989345292, 13.13725985650784, 19.87204792380568, 40.53075040663217, 16.85250252771414, 9.310383127098602, 14.34654053865356),
    @(20, 10.30247348444173, 5.138400397982204, 13.23517980964301, 20.24955283636154, 40.62157892259947, 16.84966234239457, 9.368349510365686, 14.41228102644949),
    @(21, 10.60269417340983, 5.30457569410695, 13.56351065214099, 21.46857628470577, 40.93844417558635, 16.84356135872367, 9.563031469155472, 14.63525530061323),
    @(22, 10.79637813840559, 5.40964176470384, 13.77720279263825, 22.22866616901552, 41.15404009960929, 16.84214450192212, 9.689981706480665, 14.78230950028311),
    @(23, 10.69325923141675, 5.353894648174725, 13.66325750426014, 21.82633154458858, 41.03820814493418, 16.84266519194577, 9.622266498549664, 14.70371741160222),
    @(24, 10.29752906355308, 5.135625727740451, 13.22980399417361, 20.22900810905287, 40.61654654965531, 16.84980537914268, 9.365165971975534, 14.40866239344703),
    @(25, 9.862252400723923, 4.885532722188974, 12.7610250848208, 18.34778573295695, 40.19954894129315, 16.868505791141, 9.088111680228366, 14.09761352439564)
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
